$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 122
$ws.Range("I5").Value = 115.888885
$ws.Range("K5").Value = 115.888885
$ws.Range("M5").Value = -0.8888850000000019

$ws.Range("H15").Value = 4886.511
$ws.Range("I15").Value = 4886.511
$ws.Range("K15").Value = 14659.533
$ws.Range("M15").Value = -14490.533

$ws.Range("H18").Value = 9857.143
$ws.Range("I18").Value = 4800
$ws.Range("K18").Value = 4800
$ws.Range("M18").Value = -4516

$ws.Range("H28").Value = 1071.5652
$ws.Range("I28").Value = 1074.9286
$ws.Range("J28").Value = 1066.3334
$ws.Range("K28").Value = 1074.9286
$ws.Range("L28").Value = 1066.3334
$ws.Range("M28").Value = -589.9286
$ws.Range("N28").Value = -2036.3334

$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()

$ws.Range("H107").Value = 697.4762
$ws.Range("I107").Value = 869.8571
$ws.Range("J107").Value = 352.7143
$ws.Range("K107").Value = 869.8571
$ws.Range("L107").Value = 352.7143
$ws.Range("M107").Value = 1050.1429
$ws.Range("N107").Value = -4192.7143

$ws.Range("H116").Value = 14192.854
$ws.Range("I116").Value = 14002.839
$ws.Range("J116").Value = 14781.9
$ws.Range("K116").Value = 14002.839
$ws.Range("L116").Value = 14781.9
$ws.Range("M116").Value = -10560.839
$ws.Range("N116").Value = -21665.9

$ws.Range("H118").Value = 1146.65
$ws.Range("I118").Value = 742.5833
$ws.Range("J118").Value = 1752.75
$ws.Range("K118").Value = 2227.7499
$ws.Range("L118").Value = 5258.25
$ws.Range("M118").Value = -570.7498999999998
$ws.Range("N118").Value = -8572.25

$ws.Range("H132").Value = 6804729
$ws.Range("I132").Value = 8405254
$ws.Range("J132").Value = 2498.25
$ws.Range("K132").Value = 25215762
$ws.Range("L132").Value = 7494.75
$ws.Range("M132").Value = -25213232
$ws.Range("N132").Value = -12554.75

$ws.Range("H137").Value = 13422.385
$ws.Range("I137").Value = 1438.8667
$ws.Range("J137").Value = 29763.545
$ws.Range("K137").Value = 4316.6001
$ws.Range("L137").Value = 89290.635
$ws.Range("M137").Value = -1766.6001
$ws.Range("N137").Value = -94390.635

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 41958.047
$ws.Range("I110").Value = 45848.527
$ws.Range("K110").Value = 45848.527
$ws.Range("M110").Value = -43803.527

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 365.4
$ws.Range("I22").Value = 384.33334
$ws.Range("J22").Value = 337
$ws.Range("K22").Value = 384.33334
$ws.Range("L22").Value = 337
$ws.Range("M22").Value = -211.33334
$ws.Range("N22").Value = -683

$ws.Range("H105").Value = 2381.125
$ws.Range("I105").Value = 1872.1666
$ws.Range("K105").Value = 1872.1666
$ws.Range("M105").Value = -125.1666

$ws.Range("H107").Value = 312.9
$ws.Range("I107").Value = 338
$ws.Range("K107").Value = 338
$ws.Range("M107").Value = 1582

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 900
$ws.Range("I22").Value = 300
$ws.Range("K22").Value = 300
$ws.Range("M22").Value = 50

$ws.Range("H31").Value = 2275636
$ws.Range("J31").Value = 3531.9524
$ws.Range("L31").Value = 3531.9524
$ws.Range("N31").Value = -4121.9524

$ws.Range("H34").Value = 2275636
$ws.Range("J34").Value = 3531.9524
$ws.Range("L34").Value = 3531.9524
$ws.Range("N34").Value = -3935.9524

$ws.Range("H58").Value = 930.3333
$ws.Range("I58").Value = 833.0909
$ws.Range("K58").Value = 833.0909
$ws.Range("M58").Value = -630.0909

$ws.Range("H136").Value = 930.3333
$ws.Range("I136").Value = 833.0909
$ws.Range("K136").Value = 2499.2727
$ws.Range("M136").Value = 50.72730000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H105").Value = 9699.667
$ws.Range("J105").Value = 9699.667
$ws.Range("L105").Value = 29099.001
$ws.Range("N105").Value = -34341.001

$ws.Range("H122").Value = 881.8571
$ws.Range("I122").Value = 655.7143
$ws.Range("J122").Value = 1108
$ws.Range("K122").Value = 5901.428699999999
$ws.Range("L122").Value = 9972
$ws.Range("M122").Value = -3451.428699999999
$ws.Range("N122").Value = -14872

$ws.Range("H140").Value = 1640.6364
$ws.Range("I140").Value = 1640.6364
$ws.Range("K140").Value = 4921.9092
$ws.Range("M140").Value = 258.0907999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 40054.8
$ws.Range("J20").Value = 40054.8
$ws.Range("L20").Value = 40054.8
$ws.Range("N20").Value = -40544.8

$ws.Range("H24").Value = 23322.334
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 23322.334
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 23322.334
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -23668.334

$ws.Range("H70").Value = 4999.5
$ws.Range("I70").Value = 4999.5
$ws.Range("K70").Value = 4999.5
$ws.Range("M70").Value = -4729.5

$ws.Range("H73").Value = 4999.5
$ws.Range("I73").Value = 4999.5
$ws.Range("K73").Value = 4999.5
$ws.Range("M73").Value = -4063.5

$ws.Range("H80").Value = 6769.185
$ws.Range("I80").Value = 3987.8
$ws.Range("K80").Value = 3987.8
$ws.Range("M80").Value = -2989.8

$ws.Range("H83").Value = 6769.185
$ws.Range("I83").Value = 3987.8
$ws.Range("K83").Value = 19939
$ws.Range("M83").Value = -14947

$ws.Range("H122").Value = 25002976
$ws.Range("I122").Value = 3009.2144
$ws.Range("K122").Value = 9027.643199999999
$ws.Range("M122").Value = -6577.643199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 5000250
$ws.Range("I2").Value = 500
$ws.Range("K2").Value = 500
$ws.Range("M2").Value = -388

$ws.Range("H22").Value = 1169
$ws.Range("I22").Value = 700
$ws.Range("K22").Value = 700
$ws.Range("M22").Value = -405

$ws.Range("H27").Value = 1169
$ws.Range("I27").Value = 700
$ws.Range("K27").Value = 700
$ws.Range("M27").Value = -593

$ws.Range("H61").Value = 1015.96295
$ws.Range("I61").Value = 1016.61536
$ws.Range("K61").Value = 1016.61536
$ws.Range("M61").Value = -814.61536

$ws.Range("H68").Value = 3654.9412
$ws.Range("J68").Value = 4194.273
$ws.Range("L68").Value = 4194.273
$ws.Range("N68").Value = -5692.273

$ws.Range("H71").Value = 3654.9412
$ws.Range("J71").Value = 4194.273
$ws.Range("L71").Value = 20971.365
$ws.Range("N71").Value = -28459.365

$ws.Range("H113").Value = 1015.96295
$ws.Range("I113").Value = 1016.61536
$ws.Range("K113").Value = 1016.61536
$ws.Range("M113").Value = 1153.38464

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 140000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 140000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 140000
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -141020

$ws.Range("H123").Value = 145999.5
$ws.Range("J123").Value = 145999.5
$ws.Range("L123").Value = 145999.5
$ws.Range("N123").Value = -155799.5

$ws.Range("H126").Value = 2606.7693
$ws.Range("I126").Value = 1990.25
$ws.Range("K126").Value = 5970.75
$ws.Range("M126").Value = -3500.75

$ws.Range("H136").Value = 9853.463
$ws.Range("I136").Value = 9452.979
$ws.Range("J136").Value = 12542.429
$ws.Range("K136").Value = 28358.937
$ws.Range("L136").Value = 37627.287
$ws.Range("M136").Value = -25808.937
$ws.Range("N136").Value = -42727.287
